$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume data pulled from coinranking.com
# (values mirror the scheduled GitHub Actions data refresh)

# The data columns below are stored as plain text (some price strings even
# use two '.' separators, e.g. '29.548.90', so they can't be real numbers).
# Temporarily force Text format so COM assigns the literal string as-is
# instead of silently parsing+re-rendering it as a floating point Double.
$textGuardD = $ws.Range("D2:D51")
$textGuardD.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value2 = '29.548.90'
$ws.Cells.Item(2, 5).Value2 = '  +2.25%  '
$ws.Cells.Item(3, 4).Value2 = '1.992.18'
$ws.Cells.Item(3, 5).Value2 = '  +6.01%  '
$ws.Cells.Item(4, 4).Value2 = '1.002'
$ws.Cells.Item(4, 5).Value2 = '  -0.04%  '
$ws.Cells.Item(5, 4).Value2 = '325.54'
$ws.Cells.Item(5, 5).Value2 = '  +0.19%  '
$ws.Cells.Item(6, 4).Value2 = '1.002'
$ws.Cells.Item(6, 5).Value2 = '  +0.02%  '
$ws.Cells.Item(7, 4).Value2 = '0.4684'
$ws.Cells.Item(7, 5).Value2 = '  +1.75%  '
$ws.Cells.Item(8, 4).Value2 = '0.3948'
$ws.Cells.Item(8, 5).Value2 = '  +1.76%  '
$ws.Cells.Item(9, 2).Value2 = 'Dogecoin'
$ws.Cells.Item(9, 3).Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(9, 4).Value2 = '0.07956'
$ws.Cells.Item(9, 5).Value2 = '  +1.30%  '
$ws.Cells.Item(10, 2).Value2 = 'Polygon'
$ws.Cells.Item(10, 3).Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).Value2 = '1.002'
$ws.Cells.Item(10, 5).Value2 = '  +1.66%  '
$ws.Cells.Item(11, 2).Value2 = 'Solana'
$ws.Cells.Item(11, 3).Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(11, 4).Value2 = '23.01'
$ws.Cells.Item(11, 5).Value2 = '  +5.73%  '
$ws.Cells.Item(12, 2).Value2 = 'WrappedEther'
$ws.Cells.Item(12, 3).Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value2 = '1.992.09'
$ws.Cells.Item(12, 5).Value2 = '  +9.46%  '
$ws.Cells.Item(13, 2).Value2 = 'Chainlink'
$ws.Cells.Item(13, 3).Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(13, 4).Value2 = '7.280'
$ws.Cells.Item(13, 5).Value2 = '  +4.03%  '
$ws.Cells.Item(14, 2).Value2 = 'Polkadot'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value2 = '5.874'
$ws.Cells.Item(14, 5).Value2 = '  +4.06%  '
$ws.Cells.Item(15, 2).Value2 = 'TRON'
$ws.Cells.Item(15, 3).Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(15, 4).Value2 = '0.07153'
$ws.Cells.Item(15, 5).Value2 = '  +2.69%  '
$ws.Cells.Item(16, 2).Value2 = 'Litecoin'
$ws.Cells.Item(16, 3).Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).Value2 = '88.77'
$ws.Cells.Item(16, 5).Value2 = '  +0.88%  '
$ws.Cells.Item(17, 2).Value2 = 'BinanceUSD'
$ws.Cells.Item(17, 3).Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(17, 4).Value2 = '1.004'
$ws.Cells.Item(17, 5).Value2 = '  +0.12%  '
$ws.Cells.Item(18, 2).Value2 = 'ShibaInu'
$ws.Cells.Item(18, 3).Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value2 = '0.000009945'
$ws.Cells.Item(18, 5).Value2 = '  -0.36%  '
$ws.Cells.Item(19, 2).Value2 = 'Avalanche'
$ws.Cells.Item(19, 3).Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(19, 4).Value2 = '17.40'
$ws.Cells.Item(19, 5).Value2 = '  +2.43%  '
$ws.Cells.Item(20, 2).Value2 = 'Dai'
$ws.Cells.Item(20, 3).Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(20, 4).Value2 = '1.002'
$ws.Cells.Item(20, 5).Value2 = '  -0.02%  '
$ws.Cells.Item(21, 2).Value2 = 'WrappedBTC'
$ws.Cells.Item(21, 3).Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(21, 4).Value2 = '29.680.27'
$ws.Cells.Item(21, 5).Value2 = '  +2.77%  '
$ws.Cells.Item(22, 2).Value2 = 'Uniswap'
$ws.Cells.Item(22, 3).Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value2 = '5.537'
$ws.Cells.Item(22, 5).Value2 = '  +5.81%  '
$ws.Cells.Item(23, 2).Value2 = 'Cosmos'
$ws.Cells.Item(23, 3).Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(23, 4).Value2 = '11.29'
$ws.Cells.Item(23, 5).Value2 = '  +3.07%  '
$ws.Cells.Item(24, 2).Value2 = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(24, 3).Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(24, 4).Value2 = '2.241.44'
$ws.Cells.Item(24, 5).Value2 = '  +9.28%  '
$ws.Cells.Item(25, 4).Value2 = '2.105'
$ws.Cells.Item(25, 5).Value2 = '  +0.88%  '
$ws.Cells.Item(26, 4).Value2 = '157.98'
$ws.Cells.Item(26, 5).Value2 = '  +1.06%  '
$ws.Cells.Item(27, 4).Value2 = '19.68'
$ws.Cells.Item(27, 5).Value2 = '  +1.79%  '
$ws.Cells.Item(28, 4).Value2 = '5.964'
$ws.Cells.Item(28, 5).Value2 = '  -1.08%  '
$ws.Cells.Item(29, 4).Value2 = '120.42'
$ws.Cells.Item(29, 5).Value2 = '  +2.49%  '
$ws.Cells.Item(30, 4).Value2 = '1.962'
$ws.Cells.Item(30, 5).Value2 = '  +1.82%  '
$ws.Cells.Item(31, 4).Value2 = '0.09458'
$ws.Cells.Item(31, 5).Value2 = '  +1.10%  '
$ws.Cells.Item(32, 4).Value2 = '0.9042'
$ws.Cells.Item(32, 5).Value2 = '  +0.29%  '
$ws.Cells.Item(33, 4).Value2 = '5.258'
$ws.Cells.Item(33, 5).Value2 = '  +0.07%  '
$ws.Cells.Item(35, 4).Value2 = '3.182'
$ws.Cells.Item(35, 5).Value2 = '  -2.25%  '
$ws.Cells.Item(36, 4).Value2 = '0.05843'
$ws.Cells.Item(36, 5).Value2 = '  +1.84%  '
$ws.Cells.Item(37, 4).Value2 = '1.179'
$ws.Cells.Item(37, 5).Value2 = '  -0.38%  '
$ws.Cells.Item(38, 4).Value2 = '0.02123'
$ws.Cells.Item(38, 5).Value2 = '  +2.57%  '
$ws.Cells.Item(39, 2).Value2 = 'FraxShare'
$ws.Cells.Item(39, 3).Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(39, 4).Value2 = '7.889'
$ws.Cells.Item(39, 5).Value2 = '  +3.26%  '
$ws.Cells.Item(40, 2).Value2 = 'PEPE'
$ws.Cells.Item(40, 3).Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(40, 4).Value2 = '0.000003284'
$ws.Cells.Item(40, 5).Value2 = '  +101.69%  '
$ws.Cells.Item(41, 4).Value2 = '0.5762'
$ws.Cells.Item(41, 5).Value2 = '  +2.12%  '
$ws.Cells.Item(42, 4).Value2 = '0.1830'
$ws.Cells.Item(43, 4).Value2 = '9.819'
$ws.Cells.Item(43, 5).Value2 = '  +1.69%  '
$ws.Cells.Item(44, 4).Value2 = '12.09'
$ws.Cells.Item(44, 5).Value2 = '  +2.36%  '
$ws.Cells.Item(45, 4).Value2 = '0.5373'
$ws.Cells.Item(46, 4).Value2 = '2.688'
$ws.Cells.Item(46, 5).Value2 = '  +6.21%  '
$ws.Cells.Item(47, 4).Value2 = '2.179'
$ws.Cells.Item(47, 5).Value2 = '  -4.09%  '
$ws.Cells.Item(48, 4).Value2 = '1.868'
$ws.Cells.Item(48, 5).Value2 = '  +1.48%  '
$ws.Cells.Item(49, 4).Value2 = '0.06947'
$ws.Cells.Item(49, 5).Value2 = '  -1.38%  '
$ws.Cells.Item(50, 4).Value2 = '114.59'
$ws.Cells.Item(50, 5).Value2 = '  +1.90%  '
$ws.Cells.Item(51, 4).Value2 = '0.3086'
$ws.Cells.Item(51, 5).Value2 = '  +8.08%  '

# Restore the columns to the workbook's original (unformatted) state now
# that the literal text has been committed to each cell.
$textGuardD.ClearFormats()
